# Updates cryptos list (Price / Volume(1h) columns) to the latest scraped
# values, matching the "Updated cryptos list ... with GitHub Actions" commit.
#
# Note: several "Price" values look like plain decimals (e.g. "680.42").
# Assigning such a string straight to .Value lets Excel auto-detect it as a
# number (and can even drop a significant trailing zero, e.g. "10.90" ->
# 10.9). To keep these cells as plain text - exactly like the rest of the
# Price column - we force NumberFormat "@" (Text) on just those cells
# before writing the string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.492.55'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '3.692.19'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '680.42'
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.96'
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.495'
$ws.Range("E8").Value = '  +0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.146'
$ws.Range("E9").Value = '  +0.63%  '
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").Value = '4.315.03'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.49'
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("D15").Value = '3.696.14'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").Value = '69.468.36'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("E17").Value = '  +2.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.03'
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.48'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '473.66'
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.84'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '80.23'
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("D24").Value = '3.838.87'
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.90'
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.15'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.02'
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.60'
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.99'
$ws.Range("E34").Value = '  +1.37%  '
$ws.Range("D35").Value = '3.683.31'
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.162'
$ws.Range("E36").Value = '  +2.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.45'
$ws.Range("E37").Value = '  +3.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.21'
$ws.Range("E38").Value = '  +2.51%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.26'
$ws.Range("E40").Value = '  +0.47%  '
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0906'
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '168.60'
$ws.Range("E43").Value = '  +2.23%  '
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.95'
$ws.Range("E45").Value = '  -1.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.32'
$ws.Range("E46").Value = '  +1.61%  '
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000279'
$ws.Range("E48").Value = '  +2.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.29'
$ws.Range("E49").Value = '  -1.53%  '
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.88'
$ws.Range("E51").Value = '  +0.37%  '
